$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, E across rows 2-12
$data = @{
    2  = @(96.01000000000001, 99.08, 97.91, 98.48999999999999)
    3  = @(96.01000000000001, 98.98, 98.48999999999999, 98.73)
    4  = @(96.01000000000001, 97.54000000000001, 96.41, 96.97)
    5  = @(96.01000000000001, 98.08, 94.66, 96.34)
    6  = @(96.01000000000001, 98.11, 95.17, 96.62)
    7  = @(96.01000000000001, 97.63, 94.41, 95.98999999999999)
    8  = @(96.01000000000001, 98.3, 95.69, 96.97)
    9  = @(96.01000000000001, 98.28, 95.37, 96.8)
    10 = @(96.01000000000001, 84.39, 96.73999999999999, 90.14)
    11 = @(96.01000000000001, 91.53, 94.83, 93.15000000000001)
    12 = @(96.01000000000001, 96.19, 95.97, 96.02)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
}
